$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "28.033.65", "1.002") that must remain plain text
# rather than being auto-parsed as a number. Temporarily mark the cell as Text format,
# assign the literal string, then restore the default "Normal" style so no extra
# formatting is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.033.65"
$ws.Range("E2").Value = "  -0.34%  "
Set-TextValue $ws.Range("D3") "1.869.77"
$ws.Range("E3").Value = "  -0.81%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.37%  "
Set-TextValue $ws.Range("D5") "312.45"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -0.27%  "
Set-TextValue $ws.Range("D7") "0.5116"
$ws.Range("E7").Value = "  +1.82%  "
Set-TextValue $ws.Range("D8") "0.3876"
$ws.Range("E8").Value = "  +1.37%  "
Set-TextValue $ws.Range("D9") "0.08361"
$ws.Range("E9").Value = "  -1.80%  "
Set-TextValue $ws.Range("D10") "1.112"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  +0.29%  "
Set-TextValue $ws.Range("D12") "6.169"
$ws.Range("E12").Value = "  -1.47%  "
Set-TextValue $ws.Range("D13") "1.875.43"
$ws.Range("E13").Value = "  -0.17%  "
Set-TextValue $ws.Range("D14") "20.49"
$ws.Range("E14").Value = "  -0.64%  "
Set-TextValue $ws.Range("D15") "7.271"
$ws.Range("E15").Value = "  +0.80%  "
Set-TextValue $ws.Range("D16") "1.001"
$ws.Range("E16").Value = "  -0.37%  "
Set-TextValue $ws.Range("D17") "0.00001099"
$ws.Range("E17").Value = "  +0.21%  "
Set-TextValue $ws.Range("D18") "90.87"
$ws.Range("E18").Value = "  -0.36%  "
Set-TextValue $ws.Range("D19") "0.06641"
$ws.Range("E19").Value = "  -0.12%  "
Set-TextValue $ws.Range("D20") "17.63"
$ws.Range("E20").Value = "  -2.50%  "
Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  -0.27%  "
Set-TextValue $ws.Range("D22") "6.009"
$ws.Range("E22").Value = "  -1.47%  "
Set-TextValue $ws.Range("D23") "28.063.79"
$ws.Range("E23").Value = "  -0.37%  "
Set-TextValue $ws.Range("D24") "11.07"
$ws.Range("E24").Value = "  -1.44%  "
Set-TextValue $ws.Range("D25") "2.246"
Set-TextValue $ws.Range("D26") "2.081.84"
$ws.Range("E26").Value = "  -0.60%  "
Set-TextValue $ws.Range("D27") "2.471"
$ws.Range("E27").Value = "  -5.09%  "
Set-TextValue $ws.Range("D28") "158.13"
$ws.Range("E28").Value = "  +1.32%  "
Set-TextValue $ws.Range("D29") "20.50"
Set-TextValue $ws.Range("D30") "124.68"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  -1.39%  "
Set-TextValue $ws.Range("D33") "5.871"
$ws.Range("E33").Value = "  +4.16%  "
Set-TextValue $ws.Range("D34") "3.594"
$ws.Range("E34").Value = "  -0.49%  "
Set-TextValue $ws.Range("D35") "9.413"
$ws.Range("E35").Value = "  -2.92%  "
Set-TextValue $ws.Range("D36") "0.02431"
$ws.Range("E36").Value = "  -0.86%  "
Set-TextValue $ws.Range("D37") "0.06532"
$ws.Range("E37").Value = "  +0.12%  "
Set-TextValue $ws.Range("D38") "0.2182"
$ws.Range("E38").Value = "  +0.46%  "
Set-TextValue $ws.Range("D39") "1.199"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("E41").Value = "  +2.20%  "
Set-TextValue $ws.Range("D42") "1.222"
$ws.Range("E42").Value = "  -0.93%  "
Set-TextValue $ws.Range("D43") "11.31"
$ws.Range("E43").Value = "  -0.68%  "
Set-TextValue $ws.Range("D44") "0.6074"
$ws.Range("E44").Value = "  +0.77%  "
Set-TextValue $ws.Range("D45") "12.98"
$ws.Range("E45").Value = "  -0.97%  "
Set-TextValue $ws.Range("D46") "1.279"
$ws.Range("E46").Value = "  -1.49%  "
Set-TextValue $ws.Range("D47") "3.670"
$ws.Range("E47").Value = "  -0.32%  "
Set-TextValue $ws.Range("D48") "2.004"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  -0.32%  "
Set-TextValue $ws.Range("D50") "121.02"
$ws.Range("E50").Value = "  +0.20%  "
Set-TextValue $ws.Range("D51") "77.91"
$ws.Range("E51").Value = "  -3.22%  "
